$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 1, shifting all existing rows (and their
# values) down by one.
$ws.Rows("1:1").Insert()

# Select A4 to match the resulting view state.
$ws.Range("A4").Select()
